$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the "Residuals" row (old row 6: label + Df=27, rest blank) ---
$ws.Rows(6).Delete()

# --- 1b. Drop the now-empty placeholder cells that used to pad out the merged-look
#          header row (row 1 alternates a label cell with a blank spacer; A2 was a
#          blank spacer above the "Df" label). Fully clear (not just ClearContents)
#          so the cells disappear from the sheet rather than staying as empty/styled.
foreach ($addr in @("A1", "B1", "D1", "F1", "H1", "J1", "L1", "N1", "A2")) {
    $ws.Range($addr).Clear()
}

# --- 2. Re-label the per-variable sub-headers: "F"/"P" -> "Chisq"/"Pr(>Chisq)" ---
# (row 2 holds the Df / F / P sub-header labels for each of the 6 "Soil *" blocks)
$subHeaderCols = @("C","D","E","F","G","H","I","J","K","L","M","N")
foreach ($col in $subHeaderCols) {
    $cell = $ws.Range($col + "2")
    if ($cell.Value2 -eq "F") {
        $cell.Value = "Chisq"
    } elseif ($cell.Value2 -eq "P") {
        $cell.Value = "Pr(>Chisq)"
    }
}

# --- 3. Update the statistical results themselves (new model re-run with updated values) ---
$values = @{
    "C3" = [double]"6.7290335952831999";   "D3" = [double]"9.48559783408498E-3"
    "E3" = [double]"2.8287394894881199";   "F3" = [double]"9.2591007190058E-2"
    "G3" = [double]"2.2844272786629198";   "H3" = [double]"0.13067837197047599"
    "I3" = [double]"2.5250235851722702";   "J3" = [double]"0.112053105259267"
    "K3" = [double]"1.3416909066814799";   "L3" = [double]"0.24673599343390201"
    "M3" = [double]"2.0793095787854101";   "N3" = [double]"0.14930769706744301"

    "C4" = [double]"4.14238142962986E-2";  "D4" = [double]"0.83872204119419103"
    "E4" = [double]"1.01468215258923";     "F4" = [double]"0.31378374721258201"
    "G4" = [double]"6.6637370864232697";   "H4" = [double]"9.8394358958309603E-3"
    "I4" = [double]"0.25421283913653497";  "J4" = [double]"0.61412418566564597"
    "K4" = [double]"3.17287809743508E-2";  "L4" = [double]"0.85862422388234605"
    "M4" = [double]"8.2427425931306295E-2"; "N4" = [double]"0.77403439748311698"

    "C5" = [double]"0.13529365744604099";  "D5" = [double]"0.71300532123225702"
    "E5" = [double]"6.4839733284038997E-2"; "F5" = [double]"0.79900388628772201"
    "G5" = [double]"9.9514972414990993E-2"; "H5" = [double]"0.75241246423187402"
    "I5" = [double]"0.223783796904199";    "J5" = [double]"0.63617185632787399"
    "K5" = [double]"7.8506878451526196";   "L5" = [double]"5.0801375297479896E-3"
    "M5" = [double]"2.8831516037398801";   "N5" = [double]"8.9510673895769999E-2"
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# --- 4. Resize columns to fit the new (wider) "Pr(>Chisq)" header text ---
$ws.Columns("A").ColumnWidth = 12
$ws.Columns("B").ColumnWidth = 2.1666666666666665
$ws.Columns("C:N").ColumnWidth = 11.330729166666666

# --- 5. Update the active selection to the refreshed data block ---
$ws.Range("C3:N5").Select()
